$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text updates (volume/number and date range)
$ws.Range("A8").Value = 'Volume 30   Number  14'
$ws.Range("C9").Value = 'Report Covering the Week  4/3/2023  Through  4/9/2023'

# Numeric data updates
$ws.Range("F14").Value = 1
$ws.Range("N15").Value = -40
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 400
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 54.545454545454
$ws.Range("I16").Value = 50
$ws.Range("J16").Value = 45
$ws.Range("K16").Value = 11.111111111111
$ws.Range("L16").Value = 51.515151515151
$ws.Range("M16").Value = 13.636363636363
$ws.Range("N16").Value = -82.993197278911
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = 13.333333333333
$ws.Range("I17").Value = 58
$ws.Range("J17").Value = 70
$ws.Range("K17").Value = -17.142857142857
$ws.Range("L17").Value = 5.454545454545
$ws.Range("M17").Value = 107.142857142857
$ws.Range("N17").Value = -27.5
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -38.888888888888
$ws.Range("I18").Value = 60
$ws.Range("J18").Value = 63
$ws.Range("K18").Value = -4.761904761904
$ws.Range("L18").Value = 20
$ws.Range("M18").Value = -4.761904761904
$ws.Range("N18").Value = -90.445859872611
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -38.461538461538
$ws.Range("F19").Value = 36
$ws.Range("G19").Value = 55
$ws.Range("H19").Value = -34.545454545454
$ws.Range("I19").Value = 189
$ws.Range("J19").Value = 200
$ws.Range("K19").Value = -5.5
$ws.Range("L19").Value = 92.857142857142
$ws.Range("M19").Value = 85.294117647058
$ws.Range("N19").Value = -51.908396946564
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 250
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 38.461538461538
$ws.Range("I20").Value = 54
$ws.Range("J20").Value = 48
$ws.Range("K20").Value = 12.5
$ws.Range("L20").Value = 170
$ws.Range("M20").Value = 14.893617021276
$ws.Range("N20").Value = -94.386694386694
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 45
$ws.Range("F21").Value = 101
$ws.Range("G21").Value = 112
$ws.Range("H21").Value = -9.821428571428
$ws.Range("I21").Value = 417
$ws.Range("J21").Value = 427
$ws.Range("K21").Value = -2.341920374707
$ws.Range("L21").Value = 59.770114942528
$ws.Range("M21").Value = 46.315789473684
$ws.Range("N21").Value = -82.345469940728
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 16
$ws.Range("J23").Value = 16
$ws.Range("L23").Value = 6.666666666666
$ws.Range("M23").Value = 128.571428571429
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 103
$ws.Range("G24").Value = 100
$ws.Range("H24").Value = 3
$ws.Range("I24").Value = 333
$ws.Range("J24").Value = 359
$ws.Range("K24").Value = -7.242339832869
$ws.Range("L24").Value = 8.116883116883
$ws.Range("M24").Value = 24.253731343283
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 125
$ws.Range("F25").Value = 27
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = -12.903225806451
$ws.Range("I25").Value = 92
$ws.Range("J25").Value = 97
$ws.Range("K25").Value = -5.154639175257
$ws.Range("L25").Value = 12.195121951219
$ws.Range("M25").Value = -22.033898305084
$ws.Range("G27").Value = 5
$ws.Range("F28").Value = 1
$ws.Range("N28").Value = -78.571428571428
$ws.Range("F29").Value = 1
$ws.Range("N29").Value = -72.727272727272
$ws.Range("L30").Value = -50

# Text placeholder cells (D30, E30) -- force text type to match "0" / "***.*" markers used elsewhere in the sheet
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '***.*'
